# TW edits for styling and diagrams
#
# Repositions several shapes across slides 1-3 of the FreeRADIUS MFA /
# WorkSpaces architecture diagram. All values below are expressed in
# points (PowerPoint COM units), chosen so that, after PowerPoint's
# internal Single-precision (float32) rounding and EMU conversion
# (1 pt = 12700 EMU), they reproduce the exact target EMU offsets from
# the authoritative OOXML diff.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------
# Slide 1
# ---------------------------------------------------------------
$s1 = $p.Slides.Item(1)

# Rectangle 40: off x 2926080 -> 2971800 EMU (234.0 pt)
$s1.Shapes.Item("Rectangle 40").Left = 234.0

# Rectangle 87: off x 2926080 -> 2971800 EMU (234.0 pt)
$s1.Shapes.Item("Rectangle 87").Left = 234.0

# TextBox 19 (Application Load Balancer label): off y 2836292 -> 2834640 EMU (223.20001 pt)
$s1.Shapes.Item("TextBox 19").Top = 223.20001

# Straight Arrow Connector 112: off x 5669280 -> 5715000 EMU (450.0 pt);
# ext cx 1188720 -> 1143000 EMU (90.0 pt)
$conn112 = $s1.Shapes.Item("Straight Arrow Connector 112")
$conn112.Left = 450.0
$conn112.Width = 90.0

# Straight Arrow Connector 113: ext cx 1066800 -> 1112520 EMU (87.60001 pt)
$s1.Shapes.Item("Straight Arrow Connector 113").Width = 87.60001

# Graphic 8 (picture): off x 5212080 -> 5257800 EMU (414.0 pt)
$s1.Shapes.Item("Graphic 8").Left = 414.0

# TextBox 16 (Amazon WorkSpaces label): off y 5000754 -> 5029200 EMU (396.0 pt)
$s1.Shapes.Item("TextBox 16").Top = 396.0

# Graphic 14 (picture): off x 5212080 -> 5257800 EMU (414.0 pt)
$s1.Shapes.Item("Graphic 14").Left = 414.0

# ---------------------------------------------------------------
# Slide 2
# ---------------------------------------------------------------
$s2 = $p.Slides.Item(2)

# TextBox 19 (Application Load Balancer label): off x 4475687 -> 4526280 EMU (356.40001 pt)
$s2.Shapes.Item("TextBox 19").Left = 356.40001

# TextBox 18 (AWS Managed Microsoft AD label): off x 4517639 -> 4526280 EMU (356.40001 pt)
$s2.Shapes.Item("TextBox 18").Left = 356.40001

# ---------------------------------------------------------------
# Slide 3
# ---------------------------------------------------------------
$s3 = $p.Slides.Item(3)

# TextBox 19 (Application Load Balancer label): off y 2080463 -> 2011680 EMU (158.40001 pt)
$s3.Shapes.Item("TextBox 19").Top = 158.40001
